# Doing Updates for Financials
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MSFT")

# Row 20: Non Recurring
$ws.Range("E20").Value = 5375000
$ws.Range("F20").Value = 2369000

# Row 21: Others
$ws.Range("E21").Value = 43123000
$ws.Range("F21").Value = 34747000

# Row 22: Total Operating Expenses
$ws.Range("E22").Value = 4444000
$ws.Range("F22").Value = 2486000

# Row 32: Other Items (negative mirror of row 20)
$ws.Range("E32").Value = -5375000
$ws.Range("F32").Value = -2369000

# Row 48: Property Plant and Equipment
$ws.Range("D48").Value = 36146000
$ws.Range("E48").Value = 78202000

# Row 49: Goodwill
$ws.Range("D49").Value = 43736000
$ws.Range("E49").Value = 45228000
